# Apply updated power-flow results for the 380 kV case
# (slack bus voltage setpoint B column changed from 1.05 to 1.02 pu,
# all dependent bus voltage magnitudes recomputed accordingly)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)  # B, C, D, E, F, I, J, K, L, M, N

$data = @(
    @(1.02, 1.03312892454576, 1.053403162710365, 1.043691578821516, 1.057557636787626, 1.041115032753239, 1.038255348293803, 1.056149227842278, 1.046464730418206, 1.060292279306937, 1.016723836333503),
    @(1.02, 1.033915492423133, 1.054018309200397, 1.044381608696078, 1.058307870091975, 1.041276446949004, 1.038685327062235, 1.056577919120915, 1.046966163192721, 1.060856534037247, 1.0168676811827),
    @(1.02, 1.034425263961879, 1.054416904718106, 1.044829187830686, 1.05879441549422, 1.04138004983428, 1.038963677686767, 1.05685517010841, 1.047291027575408, 1.06122207002615, 1.016960764962218),
    @(1.02, 1.034639763826675, 1.054584604947711, 1.045017607486784, 1.058999218927669, 1.041423401893196, 1.039080725018856, 1.056971691096705, 1.047427695679437, 1.061375841399514, 1.01699989846019),
    @(1.02, 1.034675790496542, 1.054612770155849, 1.045049259010008, 1.059033621492137, 1.04143066899109, 1.039100379426647, 1.056991253365618, 1.047450648387701, 1.061401666103695, 1.017006469195614),
    @(1.02, 1.034428129367865, 1.054419145025967, 1.044831704494567, 1.058797151069902, 1.041380629903673, 1.038965241567744, 1.056856727208364, 1.047292853369523, 1.06122412433594, 1.016961287862317),
    @(1.02, 1.033394580086572, 1.05361093825069, 1.043924552040364, 1.057810953726771, 1.041169757479312, 1.038400634816906, 1.056294134295054, 1.046634107412459, 1.060482882655613, 1.016772447724754),
    @(1.02, 1.031579618020689, 1.05219111508083, 1.042334423459244, 1.056081627498193, 1.040791753392567, 1.037406746964421, 1.055301764482584, 1.04547647505701, 1.059180058564732, 1.016439759031994),
    @(1.02, 1.030373981490977, 1.051247610227018, 1.041280100998536, 1.054934577905325, 1.040535482487465, 1.036744922076971, 1.054639601062784, 1.044706939865131, 1.058313862089524, 1.016218044086723),
    @(1.02, 1.029852980509666, 1.050839810872523, 1.040824960678379, 1.054439304981239, 1.040423512219301, 1.036458543539568, 1.054352758739213, 1.044374268982509, 1.057939370524464, 1.01612206365745),
    @(1.02, 1.029659616882456, 1.05068845006097, 1.040656111993615, 1.054255552284763, 1.040381771562689, 1.036352200628064, 1.054246196156025, 1.044250783391354, 1.057800356276706, 1.016086416302042),
    @(1.02, 1.029701086795378, 1.050720912274023, 1.040692321041188, 1.05429495815351, 1.040390731845995, 1.036375010126367, 1.054269054915558, 1.044277267679513, 1.057830171265653, 1.016094062589462),
    @(1.02, 1.029836993750276, 1.050827297001951, 1.040810999283031, 1.05442411155283, 1.040420064977064, 1.036449752564459, 1.054343950571468, 1.044364059918689, 1.057927877745273, 1.016119116950414),
    @(1.02, 1.029920751682341, 1.050892859318568, 1.040884148847361, 1.054503715613883, 1.040438118246699, 1.036495807973713, 1.054390094109805, 1.044417546553157, 1.057988089727584, 1.01613455432074),
    @(1.02, 1.030408580847845, 1.051274690398261, 1.041310336593138, 1.054967477382284, 1.04054289251727, 1.036763932324818, 1.054658635429334, 1.044729029713053, 1.058338728182661, 1.016224414532006),
    @(1.02, 1.030714865038201, 1.051514403801611, 1.041578046084453, 1.055258761248644, 1.040608346620339, 1.036932173097707, 1.054827053012378, 1.044924561370856, 1.058558830283536, 1.016280788143592),
    @(1.02, 1.030893616371332, 1.051654296259237, 1.041734330306381, 1.055428797811188, 1.040646428034886, 1.037030323913327, 1.05492527633619, 1.045038663955573, 1.058687267683503, 1.016313672161306),
    @(1.02, 1.030954583019744, 1.05170200806488, 1.041787641818269, 1.05548679881643, 1.040659396359769, 1.037063793955064, 1.054958765899528, 1.045077578806124, 1.058731070886142, 1.016324885126826),
    @(1.02, 1.030681993187793, 1.051488677391634, 1.041549309534225, 1.055227495202845, 1.040601334026954, 1.036914120498479, 1.054808984608653, 1.044903577255568, 1.058535209649944, 1.016274739552865),
    @(1.02, 1.029796968125392, 1.050795966186006, 1.04077604566552, 1.054386073177506, 1.040411431232852, 1.036427741917252, 1.054321896111886, 1.04433849945815, 1.057899103161648, 1.016111738949752),
    @(1.02, 1.029241440569154, 1.050361091951297, 1.040291084047506, 1.053858275413766, 1.0402911648609, 1.036122115444867, 1.054015549166069, 1.043983694711616, 1.057499671266796, 1.016009277679184),
    @(1.02, 1.029535848113911, 1.050591563649626, 1.040548055020221, 1.054137952906281, 1.040355002243637, 1.036284116388195, 1.054177957944733, 1.044171737184515, 1.057711368373131, 1.016063591925083),
    @(1.02, 1.03069684625039, 1.051500301827726, 1.041562293931179, 1.055241622564411, 1.04060450301543, 1.036922277633043, 1.054817148979242, 1.044913058912709, 1.058545882626083, 1.016277472643791),
    @(1.02, 1.032048072889199, 1.052557646078456, 1.042744502945143, 1.056527681158803, 1.040890232341439, 1.037663561483871, 1.055558425670628, 1.045775366380598, 1.059516463028282, 1.016525755512047)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $col = $cols[$j]
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
